# Generate Report for Handoff
#
# The localization-status report is regenerated. The two tracked source
# files ("9ecbf534-967f-469a-a556-bd3ccb7a1302.md" and
# "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.md") swap which report row
# describes them, and the file that now lands on row 3 ("9ecbf534...")
# has progressed from "In Translation" to "Ready for handoff" with a
# fresh handoff file/timestamp. Row 2 keeps the previously-unchanged
# "In Translation" data, now simply re-labelled for "f72f0004...".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Cells.Item(2, 1).Value = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.md"
$ws.Cells.Item(3, 1).Value = "9ecbf534-967f-469a-a556-bd3ccb7a1302.md"
$ws.Cells.Item(3, 2).Value = "Ready for handoff"
$ws.Cells.Item(3, 3).Value = "Ready for handoff"
$ws.Cells.Item(3, 4).Value = "2016-14-17 16:14:13"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "9ecbf534-967f-469a-a556-bd3ccb7a1302.md"
    }
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Cells.Item(2, 1).Value = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.md"
$ws.Cells.Item(2, 4).Value = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.c63ba4d653f502e1ea94e7d89036e70085a7864e.zh-cn.xlf"

$ws.Cells.Item(3, 1).Value = "9ecbf534-967f-469a-a556-bd3ccb7a1302.md"
$ws.Cells.Item(3, 3).Value = "Ready for handoff"
$ws.Cells.Item(3, 4).Value = "9ecbf534-967f-469a-a556-bd3ccb7a1302.7a69ee2c734182596c58358064cb221b2f4fb202.zh-cn.xlf"
$ws.Cells.Item(3, 5).Value = "2016-03-17 16:14:09"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.md"
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.c63ba4d653f502e1ea94e7d89036e70085a7864e.zh-cn.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "9ecbf534-967f-469a-a556-bd3ccb7a1302.md"
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = "9ecbf534-967f-469a-a556-bd3ccb7a1302.7a69ee2c734182596c58358064cb221b2f4fb202.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Cells.Item(2, 1).Value = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.md"
$ws.Cells.Item(2, 4).Value = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.c63ba4d653f502e1ea94e7d89036e70085a7864e.de-de.xlf"

$ws.Cells.Item(3, 1).Value = "9ecbf534-967f-469a-a556-bd3ccb7a1302.md"
$ws.Cells.Item(3, 3).Value = "Ready for handoff"
$ws.Cells.Item(3, 4).Value = "9ecbf534-967f-469a-a556-bd3ccb7a1302.7a69ee2c734182596c58358064cb221b2f4fb202.de-de.xlf"
$ws.Cells.Item(3, 5).Value = "2016-03-17 16:14:13"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.md"
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = "f72f0004-ed5e-449a-8f1e-64b6aa0789c5.c63ba4d653f502e1ea94e7d89036e70085a7864e.de-de.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "9ecbf534-967f-469a-a556-bd3ccb7a1302.md"
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = "9ecbf534-967f-469a-a556-bd3ccb7a1302.7a69ee2c734182596c58358064cb221b2f4fb202.de-de.xlf"
    }
}
